$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 375 (shifts old rows 375..461 down to 376..462)
$ws.Rows.Item(375).Insert()

# Populate the new row 375 with the new record
$ws.Cells.Item(375, 1).Value = 5
$ws.Cells.Item(375, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(375, 3).Value = "Maule"
$ws.Cells.Item(375, 4).Value = 44782
$ws.Cells.Item(375, 5).Value = 7
$ws.Cells.Item(375, 6).Value = 100112043
$ws.Cells.Item(375, 7).Value = "Pepino ensalada"
$ws.Cells.Item(375, 8).Value = "Sin especificar"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 300
$ws.Cells.Item(375, 11).Value = 20000
$ws.Cells.Item(375, 12).Value = 20000
$ws.Cells.Item(375, 13).Value = 20000
$ws.Cells.Item(375, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(375, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(375, 16).Value = 333
$ws.Cells.Item(375, 17).Value = 60
$ws.Cells.Item(375, 18).Value = "Hortaliza"
